# Remove names from the user-study results file.
# - The "Name" column header (B1) becomes "ID"
# - The participant's full name " Mike Collins" (B3) is anonymized to "M. C."
# - Column B is narrowed now that full names are no longer shown
# - The current selection on Sheet1 is reset to B2

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the name in row 3 first, then the header in row 1, so that the
# shared-strings table ends up with "M. C." before "ID" (matching the
# order in which the workbook originally listed "Mike Collins" before
# "Name" was dropped).
$ws.Range("B3").Value = "M. C."
$ws.Range("B1").Value = "ID"

# Narrow column B now that it only needs to fit short IDs/initials instead
# of full names.
$ws.Columns.Item(2).ColumnWidth = 13.666666666666666

# Reset the active selection away from the old stray M1:M1048576 selection.
$ws.Range("B2").Select()
